$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("Z1")

$ws.Range("D2").Value = "41.773.83"
$ws.Range("E2").Value = "  +2.33%  "

$ws.Range("D3").Value = "2.228.25"
$ws.Range("E3").Value = "  +0.29%  "

$ws.Range("E4").Value = "  +0.06%  "

$scratch.NumberFormat = "@"
$scratch.Value = "232.43"
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E5").Value = "  +1.58%  "

$scratch.NumberFormat = "@"
$scratch.Value = "0.622"
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E6").Value = "  -1.75%  "

$scratch.NumberFormat = "@"
$scratch.Value = "60.55"
$scratch.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E7").Value = "  -6.76%  "

$ws.Range("E8").Value = "  +0.03%  "

$scratch.NumberFormat = "@"
$scratch.Value = "0.404"
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E9").Value = "  -0.54%  "

$scratch.NumberFormat = "@"
$scratch.Value = "58.31"
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E10").Value = "  -1.48%  "

$scratch.NumberFormat = "@"
$scratch.Value = "0.0902"
$scratch.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E11").Value = "  +3.33%  "

$ws.Range("E12").Value = "  -0.44%  "

$ws.Range("D13").Value = "2.560.21"
$ws.Range("E13").Value = "  +0.39%  "

$scratch.NumberFormat = "@"
$scratch.Value = "15.61"
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E14").Value = "  -3.47%  "

$scratch.NumberFormat = "@"
$scratch.Value = "22.84"
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E15").Value = "  +2.22%  "

$scratch.NumberFormat = "@"
$scratch.Value = "0.800"
$scratch.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E16").Value = "  -3.11%  "

$scratch.NumberFormat = "@"
$scratch.Value = "5.60"
$scratch.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E17").Value = "  -0.56%  "

$ws.Range("D18").Value = "2.241.57"
$ws.Range("E18").Value = "  +1.10%  "

$ws.Range("D19").Value = "41.715.54"
$ws.Range("E19").Value = "  +2.39%  "

$ws.Range("D20").Value = "0.0₃0905"
$ws.Range("E20").Value = "  +0.34%  "

$scratch.NumberFormat = "@"
$scratch.Value = "72.35"
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E21").Value = "  -2.19%  "

$scratch.NumberFormat = "@"
$scratch.Value = "6.10"
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E22").Value = "  -1.05%  "

$scratch.NumberFormat = "@"
$scratch.Value = "247.30"
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E23").Value = "  -2.30%  "

$scratch.NumberFormat = "@"
$scratch.Value = "1.00"
$scratch.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E24").Value = "  -0.14%  "

$ws.Range("E25").Value = "  -0.31%  "

$ws.Range("E26").Value = "  -0.22%  "

$scratch.NumberFormat = "@"
$scratch.Value = "9.66"
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E27").Value = "  -0.97%  "

$scratch.NumberFormat = "@"
$scratch.Value = "168.89"
$scratch.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E28").Value = "  -2.48%  "

$scratch.NumberFormat = "@"
$scratch.Value = "0.141"
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E29").Value = "  -2.72%  "

$scratch.NumberFormat = "@"
$scratch.Value = "19.89"
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E30").Value = "  -2.16%  "

$scratch.NumberFormat = "@"
$scratch.Value = "1.39"
$scratch.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E31").Value = "  -3.27%  "

$ws.Range("E32").Value = "  -6.97%  "

$scratch.NumberFormat = "@"
$scratch.Value = "0.121"
$scratch.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E33").Value = "  -1.96%  "

$scratch.NumberFormat = "@"
$scratch.Value = "5.01"
$scratch.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E34").Value = "  +4.40%  "

$scratch.NumberFormat = "@"
$scratch.Value = "4.68"
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E35").Value = "  +0.05%  "

$scratch.NumberFormat = "@"
$scratch.Value = "0.0653"
$scratch.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E36").Value = "  +3.20%  "

$scratch.NumberFormat = "@"
$scratch.Value = "6.55"
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E37").Value = "  -8.78%  "

$scratch.NumberFormat = "@"
$scratch.Value = "2.38"
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E38").Value = "  -3.92%  "

$scratch.NumberFormat = "@"
$scratch.Value = "3.60"
$scratch.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E39").Value = "  -5.55%  "

$ws.Range("B40").Value = "TerraClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$scratch.NumberFormat = "@"
$scratch.Value = "0.000239"
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E40").Value = "  +13.60%  "

$ws.Range("B41").Value = "BinanceUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$scratch.NumberFormat = "@"
$scratch.Value = "1.00"
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E41").Value = "  +0.09%  "

$scratch.NumberFormat = "@"
$scratch.Value = "0.0240"
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E42").Value = "  +2.30%  "

$scratch.NumberFormat = "@"
$scratch.Value = "8.58"
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E43").Value = "  -1.36%  "

$ws.Range("E44").Value = "  -1.99%  "

$scratch.NumberFormat = "@"
$scratch.Value = "4.49"
$scratch.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E45").Value = "  -8.51%  "

$scratch.NumberFormat = "@"
$scratch.Value = "98.35"
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E46").Value = "  -3.53%  "

$scratch.NumberFormat = "@"
$scratch.Value = "0.0956"
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E47").Value = "  +1.52%  "

$ws.Range("D48").Value = "1.470.40"
$ws.Range("E48").Value = "  -2.71%  "

$scratch.NumberFormat = "@"
$scratch.Value = "16.55"
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E49").Value = "  -5.49%  "

$scratch.NumberFormat = "@"
$scratch.Value = "2.29"
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E50").Value = "  +8.28%  "

$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$scratch.NumberFormat = "@"
$scratch.Value = "2.75"
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E51").Value = "  -4.23%  "
